$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The product list shrank from 5 items (BRG0012..BRG0016) down to 2 new
# items (BRG0017, BRG0018). Remove the trailing rows that are no longer
# part of the data (old rows 4, 5, 6 -> Spidol/Gelas/RotiBoy/Kulkas leftovers).
$ws.Rows("4:6").Delete()

# Row 3's text cells (barang_kode / barang_nama) pick up the left-aligned
# numeric-row style (style index used by A3/D3/E3) - copy formatting from
# A3 over to B3:C3, matching the source workbook's row3 cell styles.
$ws.Range("A3").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)

# Row 2: BRG0012 / Sweater -> BRG0017 / Celana Pendek, new prices
$ws.Range("B2").Value = "BRG0017"

# Row 3: BRG0013 / Spidol -> BRG0018 / Kaos oblong, new prices and kategori_id
$ws.Range("C3").Value = "Kaos oblong"
$ws.Range("B3").Value = "BRG0018"

$ws.Range("C2").Value = "Celana Pendek"

$ws.Range("D2").Value = 40000
$ws.Range("E2").Value = 60000

$ws.Range("A3").Value = 2
$ws.Range("D3").Value = 50000
$ws.Range("E3").Value = 75000

$ws.Range("E3").Select()
